$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "51.823.85"); force text
# formatting before assignment so Excel does not coerce them into doubles, then
# clear the formatting again so no stray cell style is left behind.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range('D2').Value = '51.823.85'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '2.957.38'
$ws.Range('E3').Value = '  +3.94%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '353.27'
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('D6').Value = '112.41'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('D7').Value = '0.561'
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '0.632'
$ws.Range('E9').Value = '  +1.90%  '
$ws.Range('D10').Value = '39.54'
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('E11').Value = '  +4.85%  '
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('D13').Value = '20.03'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').Value = '7.90'
$ws.Range('E14').Value = '  +1.54%  '
$ws.Range('D15').Value = '3.411.70'
$ws.Range('E15').Value = '  +3.80%  '
$ws.Range('D16').Value = '2.953.02'
$ws.Range('E16').Value = '  +3.78%  '
$ws.Range('D17').Value = '0.988'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').Value = '51.913.33'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').Value = '7.67'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').Value = '3.33'
$ws.Range('E20').Value = '  -2.10%  '
$ws.Range('D21').Value = '14.43'
$ws.Range('E21').Value = '  +7.15%  '
$ws.Range('D22').Value = '0.0₃0988'
$ws.Range('E22').Value = '  +1.45%  '
$ws.Range('D23').Value = '71.38'
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('D24').Value = '269.49'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('E25').Value = '  +1.82%  '
$ws.Range('D26').Value = '0.181'
$ws.Range('E26').Value = '  +10.44%  '
$ws.Range('D27').Value = '27.21'
$ws.Range('E27').Value = '  +3.42%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = '0.114'
$ws.Range('E29').Value = '  +27.75%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '7.43'
$ws.Range('E30').Value = '  +18.17%  '
$ws.Range('D31').Value = '10.69'
$ws.Range('E31').Value = '  +1.39%  '
$ws.Range('E32').Value = '  +0.94%  '
$ws.Range('D33').Value = '37.58'
$ws.Range('E33').Value = '  -4.69%  '
$ws.Range('D34').Value = '6.20'
$ws.Range('E34').Value = '  +10.15%  '
$ws.Range('D35').Value = '52.96'
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('D36').Value = '0.0450'
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').Value = '3.35'
$ws.Range('E38').Value = '  +3.59%  '
$ws.Range('D39').Value = '18.85'
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('E40').Value = '  +1.79%  '
$ws.Range('D41').Value = '2.67'
$ws.Range('E41').Value = '  +5.56%  '
$ws.Range('D42').Value = '0.118'
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('D43').Value = '23.65'
$ws.Range('E43').Value = '  +5.72%  '
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').Value = '3.53'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').Value = '2.175.15'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '113.02'
$ws.Range('E48').Value = '  -8.08%  '
$ws.Range('D49').Value = '0.245'
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('D50').Value = '0.0342'
$ws.Range('E50').Value = '  +9.59%  '
$ws.Range('D51').Value = '0.936'
$ws.Range('E51').Value = '  -1.68%  '

$dRange.ClearFormats()
